# Bill of Materials update:
#  - add a new BOM line (row 16) for the Gigabit LAN Chip
#  - move the selection down to where the user would type next
#
# The workbook already has rows 3-15 populated; we only need to append
# the new component row right under the last existing entry (row 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New row 16: Gigabit LAN Chip -------------------------------------
$ws.Range("C16").Value = 10
$ws.Range("D16").Value = "Gigabit LAN Chip"
$ws.Range("E16").Value = "Gigabit Platform LAN Connect"
$ws.Range("F16").Value = "Intel"
$ws.Range("G16").Value = 82566

# Part number column (G) for this row is left-aligned text/number,
# matching the style used elsewhere in the sheet for manufacturer P/N.
$ws.Range("G16").HorizontalAlignment = -4131

# --- Move the active selection to below the newly added row ----------
$ws.Range("E17").Select()
